# Sprint 7,8,9: admin, planiller, storage
# Adds an admin-credentials table (role / email / password) to Hoja1 in
# columns G:I, rows 8-9, with the email cells turned into mailto hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: SuperAdmin ---------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("H8"), "mailto:sebastian.ravenna@gmail.com")
$ws.Range("H8").Value = "sebastian.ravenna@gmail.com"
$ws.Range("I8").Value = "admin1234"
$ws.Range("G8").Value = "SuperAdmin"

# --- Row 9: Admin ---------------------------------------------------------
$ws.Range("G9").Value = "Admin"
$ws.Hyperlinks.Add($ws.Range("H9"), "mailto:rchevi@hotmail.com")
$ws.Range("H9").Value = "rchevi@hotmail.com"
$ws.Range("I9").Value = "chevi1234"

# Restore the selection/active cell to where the author left off editing.
$ws.Range("G10").Select() | Out-Null
